$d = $word.ActiveDocument

# The document currently ends with:
#   ... "the full range of AF ablation complexity." (Conclusion paragraph)
#   <empty paragraph>                                (trailing paragraph before sectPr)
#
# We need to insert, right after the Conclusion paragraph and before the
# trailing empty paragraph:
#   - a new empty paragraph
#   - a new "NormalWeb"-styled paragraph with the new sentence(s)

$trailing = $d.Paragraphs($d.Paragraphs.Count)
$insertPoint = $trailing.Range
$insertPoint.Collapse(1)

# Insert two new paragraph marks before the trailing empty paragraph.
# Because we insert immediately before a plain/empty paragraph, the new
# paragraphs come out clean (no inherited character formatting).
$insertPoint.InsertParagraphBefore()
$insertPoint.InsertParagraphBefore()

# The first of the two new paragraphs (now second-to-last minus one) stays empty.
# The second new paragraph gets the NormalWeb style and the new text.
$newTextPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$newTextPara.Style = "NormalWeb"
$newTextPara.Range.Text = "Post-Affera non-PFA cases shifted substantially toward operators with increasing shares (Gaeta, Jae Lee, Sandesara, Kumar, Hollis), while operators with decreasing shares (especially Rashid and Fein) made up a much smaller proportion. This redistribution strongly supports operator-selection as the reason post-Affera RF cases are faster, rather than a true secular time trend."
